$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D/E remain stored as text,
# matching the original inline-string cell contents (avoids lossy numeric coercion).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.587.54'
$ws.Range('E2').Value = '  +2.71%  '
$ws.Range('D3').Value = '1.850.45'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('D4').Value = '1.036'
$ws.Range('E4').Value = '  +3.25%  '
$ws.Range('D5').Value = '321.58'
$ws.Range('E5').Value = '  +3.95%  '
$ws.Range('D6').Value = '1.031'
$ws.Range('E6').Value = '  +2.87%  '
$ws.Range('D7').Value = '0.4383'
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('D8').Value = '0.3752'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('D10').Value = '0.8761'
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').Value = '21.48'
$ws.Range('E11').Value = '  +2.47%  '
$ws.Range('D12').Value = '1.861.24'
$ws.Range('E12').Value = '  -7.06%  '
$ws.Range('D13').Value = '5.515'
$ws.Range('E13').Value = '  +2.90%  '
$ws.Range('D14').Value = '6.694'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '0.07202'
$ws.Range('E15').Value = '  +3.90%  '
$ws.Range('D16').Value = '82.74'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '1.037'
$ws.Range('E17').Value = '  +2.89%  '
$ws.Range('D18').Value = '0.000009031'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = '1.031'
$ws.Range('E19').Value = '  +2.75%  '
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').Value = '27.586.75'
$ws.Range('E21').Value = '  +2.62%  '
$ws.Range('D22').Value = '5.264'
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').Value = '2.075.43'
$ws.Range('E24').Value = '  -6.40%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '1.952'
$ws.Range('E25').Value = '  +4.35%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '157.81'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').Value = '18.73'
$ws.Range('E27').Value = '  +2.52%  '
$ws.Range('D28').Value = '5.304'
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('D30').Value = '116.30'
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = '0.09058'
$ws.Range('E31').Value = '  +1.36%  '
$ws.Range('D32').Value = '1.208'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('D33').Value = '0.7676'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').Value = '4.531'
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('D35').Value = '2.893'
$ws.Range('E35').Value = '  +3.13%  '
$ws.Range('D36').Value = '1.033'
$ws.Range('E36').Value = '  +2.60%  '
$ws.Range('D37').Value = '1.154'
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('D39').Value = '0.05289'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('D40').Value = '2.821'
$ws.Range('E40').Value = '  +5.82%  '
$ws.Range('D41').Value = '0.5175'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').Value = '6.738'
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('D44').Value = '8.583'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('D45').Value = '108.93'
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').Value = '1.715'
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('D48').Value = '0.4659'
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('D49').Value = '0.06394'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').Value = '1.894'
$ws.Range('E50').Value = '  +4.22%  '
$ws.Range('E51').Value = '  +5.19%  '
